# add language id for property name
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")

$lpids = @(
    "LPID_SUCKBLOOD",
    "LPID_REFLECTDAMAGE",
    "LPID_CRITICAL",
    "LPID_MAXHP",
    "LPID_MAXMP",
    "LPID_MAXSP",
    "LPID_HPREGEN",
    "LPID_SPREGEN",
    "LPID_MPREGEN",
    "LPID_ATK_VALUE",
    "LPID_DEF_VALUE",
    "LPID_MOVE_SPEED",
    "LPID_ATK_SPEED",
    "LPID_ATK_FIRE",
    "LPID_ATK_LIGHT",
    "LPID_ATK_WIND",
    "LPID_ATK_ICE",
    "LPID_ATK_POISON",
    "LPID_DEF_FIRE",
    "LPID_DEF_LIGHT",
    "LPID_DEF_WIND",
    "LPID_DEF_ICE",
    "LPID_DEF_POISON",
    "LPID_DIZZY_GATE",
    "LPID_MOVE_GATE",
    "LPID_SKILL_GATE",
    "LPID_PHYSICAL_GATE",
    "LPID_MAGIC_GATE",
    "LPID_BUFF_GATE"
)

for ($i = 0; $i -lt $lpids.Length; $i++) {
    $col = $i + 2   # column B = 2
    $ws.Cells.Item(9, $col).Value = $lpids[$i]
}

# Column B now holds longer identifier text (e.g. "LPID_REFLECTDAMAGE") instead of
# the short CJK labels it used to hold, so it no longer fits with columns C:D -
# widen it to fit its new content.
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

# Match the author's final selection state (active cell AD9)
$ws.Activate() | Out-Null
$ws.Range("AD9").Select() | Out-Null
